# Add a new data row (row 4) to Sheet1, mirroring the existing rows 2-3:
#   Equipment Number | Vessel | Voyage | WONumber | ReferenceNumber | BOLNumber
#
# D4/E4 hold a purely-numeric-looking reference number ("7032005141") that
# must be stored as text (shared string), same as the other rows, rather
# than being auto-coerced into a number -- so format those two cells as
# Text before assigning their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4:E4").NumberFormat = "@"

$ws.Range("A4").Value = "PCIU1438389"
$ws.Range("B4").Value = "KOTA PERWIRA"
$ws.Range("C4").Value = "0004E"
$ws.Range("D4").Value = "7032005141"
$ws.Range("E4").Value = "7032005141"
$ws.Range("F4").Value = "DEL900004900"
